$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet


# --- Row 6: Issuer Name ---
$ws.Range("E6").Value = "EVER FRESH CAMERON SDN. BHD."

# --- Row 7: Stage / Date (CRA Report) / Name fields ---
$ws.Range("J7").Value = "'2025-10-28"
$ws.Range("M7").Value = "EVER FRESH CAMERON SDN. BHD."
$ws.Range("O7").Value = ""

# --- Row 11: Scoring by CRA Agency (Issuer's Credit Agency Score) ---
$ws.Range("M11").Value = "'300"
$ws.Range("O11").Value = "'543"
$ws.Range("Q11").Value = "'644"
$ws.Range("S11").Value = "'575"
$ws.Range("S11").Style = "Normal"
$ws.Range("U11").Value = "'575"
$ws.Range("U11").Style = "Normal"

# --- Row 12: Scoring by CRA Agency (Credit Score Equivalent) ---
$ws.Range("M12").Value = "F"
$ws.Range("Q12").Value = "B"
$ws.Range("S12").Value = "C"
$ws.Range("S12").Style = "Normal"
$ws.Range("U12").Value = "C"
$ws.Range("U12").Style = "Normal"

# --- Row 13: Business has been in operations for at least THREE (3) years ---
$ws.Range("M13").Value = "'2013"

# --- Row 18: Credit Applications Approved for Last 12 months ---
$ws.Range("O18").Value = "'3"
$ws.Range("Q18").Value = "'3"
$ws.Range("S18").Value = "'1"
$ws.Range("S18").Style = "Normal"
$ws.Range("U18").Value = "'1"
$ws.Range("U18").Style = "Normal"

# --- Row 19: Credit Applications Pending ---
$ws.Range("Q19").Value = "'0"
$ws.Range("S19").Value = "'1"
$ws.Range("S19").Style = "Normal"
$ws.Range("U19").Value = "'1"
$ws.Range("U19").Style = "Normal"

# --- Row 20: Legal Action taken (from Banking) ---
$ws.Range("Q20").Value = "'0"
$ws.Range("S20").Value = "'0"
$ws.Range("S20").Style = "Normal"
$ws.Range("U20").Value = "'0"
$ws.Range("U20").Style = "Normal"

# --- Row 21: Existing No. of Facility (from Banking) ---
$ws.Range("M21").Value = "'35"
$ws.Range("O21").Value = "'11"
$ws.Range("Q21").Value = "'12"
$ws.Range("S21").Value = "'1"
$ws.Range("S21").Style = "Normal"
$ws.Range("U21").Value = "'1"
$ws.Range("U21").Style = "Normal"

# --- Row 22: Legal Suits ---
$ws.Range("M22").Value = "'2"
$ws.Range("O22").Value = "'0"
$ws.Range("Q22").Value = "'0"
$ws.Range("S22").Value = "'0"
$ws.Range("S22").Style = "Normal"
$ws.Range("U22").Value = "'0"
$ws.Range("U22").Style = "Normal"

# --- Row 23: Legal Case - Status ---
$ws.Range("M23").Value = "Yes, Yes, No"
$ws.Range("O23").Value = "No, No, No"
$ws.Range("Q23").Value = "No, No, No"
$ws.Range("S23").Value = "No, No, No"
$ws.Range("S23").Style = "Normal"
$ws.Range("U23").Value = "No, No, No"
$ws.Range("U23").Style = "Normal"

# --- Row 24: Trade / Credit Reference ---
$ws.Range("M24").Value = "'5"
$ws.Range("O24").Value = "'5"
$ws.Range("Q24").Value = "'5"
$ws.Range("S24").Value = "'5"
$ws.Range("S24").Style = "Normal"
$ws.Range("U24").Value = "'5"
$ws.Range("U24").Style = "Normal"

# --- Row 25: Total Enquiries for Last 12 months ---
$ws.Range("O25").Value = "'1"
$ws.Range("Q25").Value = "'2"
$ws.Range("S25").Value = "'2"
$ws.Range("S25").Style = "Normal"
$ws.Range("U25").Value = "'2"
$ws.Range("U25").Style = "Normal"

# --- Row 26: Special Attention Account ---
$ws.Range("Q26").Value = "'0"
$ws.Range("S26").Value = "'0"
$ws.Range("S26").Style = "Normal"
$ws.Range("U26").Value = "'0"
$ws.Range("U26").Style = "Normal"

# --- Row 27: Summary of Total Liabilities (Outstanding) ---
$ws.Range("M27").Value = "'30667835"
$ws.Range("O27").Value = "'1349839"
$ws.Range("Q27").Value = "'1278562"
$ws.Range("S27").Value = "'17053"
$ws.Range("S27").Style = "Normal"
$ws.Range("U27").Value = "'17053"
$ws.Range("U27").Style = "Normal"

# --- Row 28: Summary of Total Liabilities (Total Limit) ---
$ws.Range("M28").Value = "'34070376"
$ws.Range("O28").Value = "'1669243"
$ws.Range("Q28").Value = "'1707145"
$ws.Range("S28").Value = "'40000"
$ws.Range("S28").Style = "Normal"
$ws.Range("U28").Value = "'40000"
$ws.Range("U28").Style = "Normal"

# --- Row 29: Overdraft facility outstanding amount does not exceed limit ---
$ws.Range("Q29").Value = "No"
$ws.Range("S29").Value = "No"
$ws.Range("S29").Style = "Normal"
$ws.Range("U29").Value = "No"
$ws.Range("U29").Style = "Normal"

# --- Row 30: Issuer's Total Banking Outstanding Facilities does not exceed Total Banking Limit ---
$ws.Range("M30").Value = "YES, outstanding: 30667835.0, limit: 34070376.0"
$ws.Range("O30").Value = "YES, outstanding: 30667835.0, limit: 34070376.0"
$ws.Range("Q30").Value = "YES, outstanding: 30667835.0, limit: 34070376.0"
$ws.Range("S30").Value = "YES, outstanding: 30667835.0, limit: 34070376.0"
$ws.Range("S30").Style = "Normal"
$ws.Range("U30").Value = "YES, outstanding: 30667835.0, limit: 34070376.0"
$ws.Range("U30").Style = "Normal"

# --- Row 31: CCRIS Loan Account - Conduct Count ---
$ws.Range("M31").Value = "current 1 month MIA1: 34, MIA2: 0, MIA3: 1, MIA4+: 0 and /or past 6 months MIA1: 61, MIA2: 0, MIA3: 2, MIA4+: 0"
$ws.Range("O31").Value = "current 1 month MIA1: 4, MIA2: 0, MIA3: 0, MIA4+: 0 and /or past 6 months MIA1: 10, MIA2: 1, MIA3: 0, MIA4+: 0"
$ws.Range("Q31").Value = "current 1 month MIA1: 1, MIA2: 0, MIA3: 0, MIA4+: 0 and /or past 6 months MIA1: 8, MIA2: 0, MIA3: 0, MIA4+: 0"
$ws.Range("S31").Value = "current 1 month MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 0 and /or past 6 months MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 0"
$ws.Range("S31").Style = "Normal"
$ws.Range("U31").Value = "current 1 month MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 0 and /or past 6 months MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 0"
$ws.Range("U31").Style = "Normal"

# --- Row 33: Issuer's Total Non-Bank Lender Outstanding Facilities does not exceed limit ---
$ws.Range("Q33").Value = "NO"
$ws.Range("S33").Value = "NO"
$ws.Range("S33").Style = "Normal"
$ws.Range("U33").Value = "NO"
$ws.Range("U33").Style = "Normal"

# --- Row 34: Non-Bank Lender Credit Information (NLCI) - Conduct Count ---
$ws.Range("M34").Value = "current 1 month MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 2 and /or past 6 months MIA1: 2, MIA2: 4, MIA3: 2, MIA4+: 4"
$ws.Range("O34").Value = "current 1 month MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 2 and /or past 6 months MIA1: 2, MIA2: 4, MIA3: 2, MIA4+: 4"
$ws.Range("Q34").Value = "current 1 month MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 2 and /or past 6 months MIA1: 2, MIA2: 4, MIA3: 2, MIA4+: 4"
$ws.Range("S34").Value = "current 1 month MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 2 and /or past 6 months MIA1: 2, MIA2: 4, MIA3: 2, MIA4+: 4"
$ws.Range("S34").Style = "Normal"
$ws.Range("U34").Value = "current 1 month MIA1: 0, MIA2: 0, MIA3: 0, MIA4+: 2 and /or past 6 months MIA1: 2, MIA2: 4, MIA3: 2, MIA4+: 4"
$ws.Range("U34").Style = "Normal"

# --- Row 35: Non-Bank Lender Credit Information (NLCI) - Legal Status ---
$ws.Range("M35").Value = "LOD"
$ws.Range("O35").Value = "LOD"
$ws.Range("Q35").Value = "LOD"
$ws.Range("S35").Value = "LOD"
$ws.Range("S35").Style = "Normal"
$ws.Range("U35").Value = "LOD"
$ws.Range("U35").Style = "Normal"
